# Petty cash book — "Update 8-Jul-2021, midday update."
# Applies the 8-Jul-2021 transactions to the "Sheet1" (Buku Kas Umum) ledger:
#   - D20: Uang makan top-up (+260,000 on top of the usual 60,000)
#   - D21: TRANSFER BCA total grows by +1,740,000
#   - C23: A/R collections grow by +14,393,000
#   - New rows 25-28: SALES cash/retail, SELISIH lebih, SETOR KE BANK, Wages Expense (8-Jul)
#   - Move the frozen-pane scroll position / active selection down to the new data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 20: Uang makan (4) — additional 260,000 ---
$ws.Range("D20").Formula = "=60000+260000"

# --- Row 21: TRANSFER BCA — additional 1,740,000 ---
$ws.Range("D21").Formula = "=1405000+1864000+6027200+9027000+3000000+9027000+2606000+1897000+1740000"

# --- Row 23: A/R — additional 14,393,000 ---
$ws.Range("C23").Formula = "=6027200+9027000+3000000+9027000+29760000+2281500+1897000+14393000"

# --- Row 25: SALES - cash/retail ---
$ws.Range("B25").Value = "SALES - cash/retail"
$ws.Range("C25").Formula = "=41437225-24476225-14393000"

# --- Row 26: SELISIH - lebih ---
$ws.Range("B26").Value = "SELISIH - lebih"
$ws.Range("C26").Value = 10000

# --- Row 27: SETOR KE BANK ---
$ws.Range("B27").Value = "SETOR KE BANK"
$ws.Range("D27").Formula = "=41000000"

# --- Row 28: new day, 8-Jul-2021 — Uang makan (4) ---
$ws.Range("A28").Value = 44385
$ws.Range("B28").Value = "Wages Expense"
$ws.Range("D28").Formula = "=60000"

# --- Scroll the frozen pane down to the new bottom of the data and move selection ---
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select()

Write-Output "Applied 8-Jul-2021 midday update."
